$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember how many rows of data existed before the edit (header + 20
# data rows, i.e. A1:C21) so the new final row can be computed below.
$originalLastRow = $ws.UsedRange.Rows.Count

# 9 new data rows get inserted right after the header row (row 1), which
# pushes the existing data rows (old rows 2-21) down to rows 11-30. One
# more brand-new row is then appended at the very end (row 31).
$insertedRows = @(
    @(-0.188694953918457, -0.0127399563789367, 0.0153613984584808),
    @(-0.0261173248291015, -0.1474769711494445, 0.0655251443386077),
    @(-0.1960973739624023, 0.0549294650554657, 0.0360765755176544),
    @(-0.06610202789306641, -0.1787786185741424, 0.0745508223772049),
    @(0.0234136581420898, 0.0270741879940032, 0.2239813506603241),
    @(0.11651611328125, -0.4856438636779785, 0.5658785104751587),
    @(0.0557413101196289, 0.3574482798576355, 0.2321825623512268),
    @(0.3619680404663086, 0.0124948024749755, 0.3587799966335296),
    @(-0.2529764175415039, 0.1160029470920562, -0.09882223606109609)
)

$appendedRow = @(-0.8564167022705078, -0.1756476759910583, -1.401212096214294)

# Insert the new rows above the current row 2, then clear any formatting
# the insert operation may have copied down from the header row so the
# new cells stay unstyled (matching the rest of the plain data rows).
$insertEndRow = 1 + $insertedRows.Count
$ws.Range("A2:C$insertEndRow").EntireRow.Insert()
$ws.Range("A2:C$insertEndRow").ClearFormats()

for ($i = 0; $i -lt $insertedRows.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $insertedRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $insertedRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $insertedRows[$i][2]
}

# Append the final brand-new row after the (now shifted) existing data.
$lastRow = $originalLastRow + $insertedRows.Count + 1
$ws.Cells.Item($lastRow, 1).Value = $appendedRow[0]
$ws.Cells.Item($lastRow, 2).Value = $appendedRow[1]
$ws.Cells.Item($lastRow, 3).Value = $appendedRow[2]
